$wb = $excel.ActiveWorkbook

# Rename the "sheet1-19nov" sheet to "sheet1-20nov"
$ws = $wb.Worksheets.Item("sheet1-19nov")
$ws.Name = "sheet1-20nov"

# Fill in the new "rice tracking" K-column values on rows 8-17.
# K8 already has a cell (formatted); rows 9-17 need new cells matching
# the row's existing style, so copy it from the neighboring J cell.
$ws.Range("K8").Value = 453

$riceValues = @{ 9 = 197; 10 = 4539; 11 = 0; 12 = 2; 13 = 514; 14 = 273; 15 = 3; 16 = 1; 17 = 43 }
foreach ($row in $riceValues.Keys) {
    $kCell = $ws.Range("K$row")
    $kCell.Value = $riceValues[$row]
    $kCell.Style = $ws.Range("J$row").Style
}

# Update the view: active cell/selection moved from J17 to K17
# (also nudge the scrolled-to column left by one, from G to F)
$ws.Activate()
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("K17").Select()

